# Recipes_Database.xlsx edit:
#   "Now the script can add 1 or more recipes"
#
# Marks Salmon_Pasta as a "selected" recipe on the Control_sheet (adds the
# "x" marker in column B, mirroring Pesto_Pasta's existing marker), fixes
# the capitalisation of a couple of ingredient names on the Salmon_Pasta
# sheet, and normalises two mis-cased/duplicated ingredient category names
# ("Spice"/"chives") back onto their canonical shared entries
# ("Spices"/"Chives"). Finishes with the selection back on the
# Control_sheet, which becomes the active tab.

$wb = $excel.ActiveWorkbook

$controlSheet = $wb.Worksheets.Item("Control_sheet")
$salmonSheet  = $wb.Worksheets.Item("Salmon_Pasta")

# --- Salmon_Pasta: ingredient name clean-up -------------------------------
# Capitalisation fixes.
$salmonSheet.Range("A2").Value = "Wholewheat Pasta"
$salmonSheet.Range("A3").Value = "Shallot"

# Re-point the mis-cased/duplicate category entries to their canonical
# ("Spices" / "Chives") counterparts already used elsewhere in the sheet.
$salmonSheet.Range("B7").Value = "Spices"
$salmonSheet.Range("A8").Value = "Chives"

# Leave the cursor where the edits finished.
[void]$salmonSheet.Range("A3").Select()

# --- Control_sheet: mark Salmon_Pasta as selected -------------------------
$controlSheet.Range("B4").Value = "x"

# Control_sheet becomes the active sheet/tab, cursor below the table.
[void]$controlSheet.Activate()
[void]$controlSheet.Range("B6").Select()
